$wb = $excel.ActiveWorkbook

$rb = $wb.Worksheets.Item("RB")
$wr = $wb.Worksheets.Item("WR")

# --- WR sheet: insert a new row for "K.Yeboah" above the existing row 11,
#     pushing the former row 11 ("K.Yeboah") down to row 12, then rename
#     that pushed-down row's player to "T.Black" ---
$wr.Rows("11").Insert()

$wr.Range("A11").Value = "K.Yeboah"
$wr.Range("B11:J11").Value = 0

$wr.Range("A12").Value = "T.Black"

# --- RB sheet: rename "T.Johnson" -> "Ty.Johnson" ---
$rb.Range("A3").Value = "Ty.Johnson"

# Selections left by the editing session
$wr.Range("J13").Select()
$rb.Range("A4").Select()

# --- Activate RB tab (it was the last sheet worked on / viewed) ---
$rb.Activate()
$rb.Range("A4").Select()
